$d = $word.ActiveDocument

function FindReplace([string]$old, [string]$new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

# 1. Date change
FindReplace "2021-08-05" "2021-08-06"

# 2. Executive summary bullet 1 - add "since 1978"
FindReplace "The recent African Swine Fever (ASF) outbreak in Dominican Republic (D.R.) is the first detection of the virus in the Americas." "The recent African Swine Fever (ASF) outbreak in Dominican Republic (D.R.) is the first detection of the virus in the Americas since 1978."

# 3. Executive summary bullet 3 - reword
FindReplace "The D.R. was average among countries in the Americas for risk of ASF introduction. Many other countries in the Americas including the U.S. are more susceptible to introduction." "Prior to the current introduction, the D.R. was average among countries in the Americas for risk of ASF introduction. Many other countries in the Americas, including the U.S., were calculated to be more susceptible to introduction."

# 4. Split the final executive-summary bullet into three bullets, with "among" bolded
$targetText = "The D.R. is now among the most likely sources of ASF spread to elsewhere in the Americas. However, this risk does not greatly exceed other potential source countries from Africa and Asia (e.g., South Africa, Philippines, China, and Russia). Risk mitigation strategies should incorporate source control from the D.R. but continue to focus equally on all major global sources."
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq $targetText) {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not locate target paragraph for split"
}

$fragBody = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1001"/></w:numPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t xml:space="preserve">The D.R. is now</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">among</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">the most likely sources of ASF spread to elsewhere in the Americas. This emergent risk warrants increased biosecurity measures, such as increasing inspections of flights and boats from the D.R.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1001"/></w:numPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t xml:space="preserve">However, this risk does not greatly exceed other potential source countries from Africa and Asia (e.g., South Africa, Philippines, China, and Russia).</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1001"/></w:numPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t xml:space="preserve">Risk mitigation strategies should incorporate source control from the D.R. but continue to focus equally on other identified major global sources. The introduction to the D.R. reaffirms the risk from other, more distant sources.</w:t></w:r></w:p>'

$xmlPkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $fragBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetPara.Range.InsertXML($xmlPkg)

# 5. Introduction paragraph - add "in over 40 years"
FindReplace "In July 2021, the U.S. Department of Agriculture confirmed an outbreak of African Swine Fever (ASF) in pigs in the Dominican Republic (D.R.), marking the first detection of the virus in the western hemisphere and raising concerns about potential transmission to the United States. The source of the outbreak in the D.R. is currently unknown." "In July 2021, the U.S. Department of Agriculture confirmed an outbreak of African Swine Fever (ASF) in pigs in the Dominican Republic (D.R.), marking the first detection of the virus in the western hemisphere in over 40 years and raising concerns about potential transmission to the United States. The source of the outbreak in the D.R. is currently unknown."

# 6. Remove stray "other" before "countries in the Americas are subject"
FindReplace "), has the highest estimated risk in the Americas. The D.R. was a middle-risk country prior to July 2021. This indicates that arrival of ASF was likely not due to unique D.R. vulnerabilities but that most other countries in the Americas are subject to the same import risks." "), has the highest estimated risk in the Americas. The D.R. was a middle-risk country prior to July 2021. This indicates that arrival of ASF was likely not due to unique D.R. vulnerabilities but that most countries in the Americas are subject to the same import risks."

# 9. Add "(Figure 3)" reference
FindReplace "of risk. The D.R. is is now among the most likely origin countries for potential ASF spread to North American countries (USA, Mexico, Canada) and higher-risk countries in the Caribbean and South America (Haiti, Brazil, Colombia, Cuba). This emergent risk warrants increased biosecurity measures, such as increasing inspections of flights and boats from the D.R.." "of risk. The D.R. is is now among the most likely origin countries for potential ASF spread to North American countries (USA, Mexico, Canada) and higher-risk countries in the Caribbean and South America (Haiti, Brazil, Colombia, Cuba) (Figure 3). This emergent risk warrants increased biosecurity measures, such as increasing inspections of flights and boats from the D.R.."

# 11. Fix typo "contries" -> "countries"
FindReplace "REPEL estimates of relative contribution of ASF source contries to ASF import probability to select countries in the Americas" "REPEL estimates of relative contribution of ASF source countries to ASF import probability to select countries in the Americas"

Write-Output "All edits applied"
